$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.334.50'
Set-TextValue $ws.Range("E2") '  -4.42%  '

Set-TextValue $ws.Range("D3") '1.762.37'
Set-TextValue $ws.Range("E3") '  -3.94%  '

Set-TextValue $ws.Range("E4") '  -0.09%  '

Set-TextValue $ws.Range("E5") '  -0.05%  '

Set-TextValue $ws.Range("D6") '304.09'
Set-TextValue $ws.Range("E6") '  -2.55%  '

Set-TextValue $ws.Range("D7") '0.4257'
Set-TextValue $ws.Range("E7") '  -0.44%  '

Set-TextValue $ws.Range("D8") '0.3618'
Set-TextValue $ws.Range("E8") '  -1.01%  '

Set-TextValue $ws.Range("D9") '0.07054'
Set-TextValue $ws.Range("E9") '  -2.99%  '

Set-TextValue $ws.Range("D10") '0.8302'
Set-TextValue $ws.Range("E10") '  -3.95%  '

Set-TextValue $ws.Range("E11") '  -2.56%  '

Set-TextValue $ws.Range("D12") '1.750.86'
Set-TextValue $ws.Range("E12") '  +0.52%  '

Set-TextValue $ws.Range("D13") '5.233'
Set-TextValue $ws.Range("E13") '  -4.31%  '

Set-TextValue $ws.Range("D14") '6.390'
Set-TextValue $ws.Range("E14") '  -2.01%  '

Set-TextValue $ws.Range("D15") '0.06777'
Set-TextValue $ws.Range("E15") '  -2.76%  '

Set-TextValue $ws.Range("D16") '1.004'
Set-TextValue $ws.Range("E16") '  +0.07%  '

Set-TextValue $ws.Range("D17") '79.13'
Set-TextValue $ws.Range("E17") '  -1.91%  '

Set-TextValue $ws.Range("D18") '0.000008614'
Set-TextValue $ws.Range("E18") '  -3.41%  '

Set-TextValue $ws.Range("E19") '  +0.04%  '

Set-TextValue $ws.Range("D20") '14.93'
Set-TextValue $ws.Range("E20") '  -3.09%  '

Set-TextValue $ws.Range("D21") '25.800.07'
Set-TextValue $ws.Range("E21") '  -5.53%  '

Set-TextValue $ws.Range("D22") '4.995'
Set-TextValue $ws.Range("E22") '  -3.06%  '

Set-TextValue $ws.Range("D23") '11.07'
Set-TextValue $ws.Range("E23") '  +1.60%  '

Set-TextValue $ws.Range("D24") '1.930.46'
Set-TextValue $ws.Range("E24") '  -2.74%  '

Set-TextValue $ws.Range("D25") '1.904'
Set-TextValue $ws.Range("E25") '  -4.38%  '

Set-TextValue $ws.Range("D26") '152.01'
Set-TextValue $ws.Range("E26") '  -1.92%  '

Set-TextValue $ws.Range("D27") '18.11'
Set-TextValue $ws.Range("E27") '  -4.04%  '

Set-TextValue $ws.Range("D28") '114.68'
Set-TextValue $ws.Range("E28") '  +0.40%  '

Set-TextValue $ws.Range("D29") '5.007'
Set-TextValue $ws.Range("E29") '  -2.84%  '

Set-TextValue $ws.Range("D30") '1.672'
Set-TextValue $ws.Range("E30") '  -8.01%  '

Set-TextValue $ws.Range("D31") '0.08876'
Set-TextValue $ws.Range("E31") '  +0.26%  '

Set-TextValue $ws.Range("D32") '0.7201'
Set-TextValue $ws.Range("E32") '  -3.82%  '

Set-TextValue $ws.Range("D33") '1.116'
Set-TextValue $ws.Range("E33") '  -1.46%  '

Set-TextValue $ws.Range("E34") '  -5.39%  '

Set-TextValue $ws.Range("D35") '0.9994'
Set-TextValue $ws.Range("E35") '  -0.18%  '

Set-TextValue $ws.Range("D36") '2.707'
Set-TextValue $ws.Range("E36") '  -9.48%  '

Set-TextValue $ws.Range("D37") '1.068'
Set-TextValue $ws.Range("E37") '  -2.72%  '

Set-TextValue $ws.Range("D38") '0.05088'
Set-TextValue $ws.Range("E38") '  -4.42%  '

Set-TextValue $ws.Range("D39") '0.01881'
Set-TextValue $ws.Range("E39") '  -2.87%  '

Set-TextValue $ws.Range("B40") 'TheSandbox'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D40") '0.4891'
Set-TextValue $ws.Range("E40") '  -3.50%  '

Set-TextValue $ws.Range("B41") 'Algorand'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D41") '0.1599'
Set-TextValue $ws.Range("E41") '  -3.05%  '

Set-TextValue $ws.Range("D42") '6.175'
Set-TextValue $ws.Range("E42") '  -4.27%  '

Set-TextValue $ws.Range("D43") '2.466'
Set-TextValue $ws.Range("E43") '  -11.98%  '

Set-TextValue $ws.Range("D44") '7.971'
Set-TextValue $ws.Range("E44") '  -4.31%  '

Set-TextValue $ws.Range("D45") '104.51'
Set-TextValue $ws.Range("E45") '  -0.86%  '

Set-TextValue $ws.Range("D46") '1.000'
Set-TextValue $ws.Range("E46") '  +0.02%  '

Set-TextValue $ws.Range("D47") '10.00'
Set-TextValue $ws.Range("E47") '  -3.83%  '

Set-TextValue $ws.Range("D48") '0.06180'
Set-TextValue $ws.Range("E48") '  -4.54%  '

Set-TextValue $ws.Range("D49") '0.4453'
Set-TextValue $ws.Range("E49") '  -5.05%  '

Set-TextValue $ws.Range("E50") '  -3.46%  '

Set-TextValue $ws.Range("D51") '1.709'
Set-TextValue $ws.Range("E51") '  -1.80%  '
